# Auto-generated edit script applying the cryptos.xlsx diff
# (GitHub Actions "Updated cryptos list" style update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values look numeric ("29.852.25", "0.07202", "1.000", ...
# with thousands separators / significant trailing zeros) and must stay
# literal text, exactly as authored, so force Text format before writing.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.852.25"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.893.58"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7897"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.46"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3149"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.32"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07202"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08089"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7650"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.499"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.880.57"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.48"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.142"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.861.83"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.93"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.83"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007782"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.169.23"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.111"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1646"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.399"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.32"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.72"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.051"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.404"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.505"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05566"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7427"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.000"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.612"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01922"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.774"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.141.69"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.86"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4422"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.874"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8508"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "104.14"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.06"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.879"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.460"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.014"

# Coin name / link / volume columns are plain text already, no coercion risk.
$ws.Range("E2").Value = "  -1.13%  "
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("E5").Value = "  -4.68%  "
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  -3.72%  "
$ws.Range("E9").Value = "  -5.72%  "
$ws.Range("E10").Value = "  +2.01%  "
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("E13").Value = "  +3.97%  "
$ws.Range("E14").Value = "  -1.61%  "
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("E16").Value = "  +3.91%  "
$ws.Range("E17").Value = "  -1.13%  "
$ws.Range("E18").Value = "  -1.99%  "
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("E23").Value = "  +14.96%  "
$ws.Range("E25").Value = "  -2.42%  "
$ws.Range("E26").Value = "  +0.84%  "
$ws.Range("E27").Value = "  -1.98%  "
$ws.Range("E28").Value = "  -1.54%  "
$ws.Range("E29").Value = "  -2.49%  "
$ws.Range("E30").Value = "  +2.55%  "
$ws.Range("E31").Value = "  +1.72%  "
$ws.Range("E32").Value = "  +4.59%  "
$ws.Range("E33").Value = "  +0.48%  "
$ws.Range("E34").Value = "  -8.63%  "
$ws.Range("E35").Value = "  -0.25%  "
$ws.Range("E36").Value = "  +1.07%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("E38").Value = "  -4.00%  "
$ws.Range("E39").Value = "  -0.71%  "
$ws.Range("E40").Value = "  -0.79%  "
$ws.Range("E41").Value = "  +13.01%  "
$ws.Range("E42").Value = "  +0.74%  "
$ws.Range("E43").Value = "  -0.97%  "
$ws.Range("E44").Value = "  -1.47%  "
$ws.Range("E45").Value = "  -0.41%  "
$ws.Range("E46").Value = "  +1.58%  "
$ws.Range("E47").Value = "  -0.10%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E48").Value = "  +1.93%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E49").Value = "  -1.44%  "
$ws.Range("E50").Value = "  -1.75%  "
$ws.Range("E51").Value = "  +10.05%  "
